$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8826830387115479
$ws.Range("B1").Value = 2.983670949935913
$ws.Range("C1").Value = 4.459150791168213
$ws.Range("D1").Value = 3.021851778030396
$ws.Range("E1").Value = 1.412309408187866
